$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StateCounters")

$ws.Range("B2").Value = 9423
$ws.Range("C2").Value = 2275
